$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column G (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 20 }

# Add the new "Save" header in H1, matching the header style used by the
# other header cells (bold font + thin border, centered/top aligned).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Find the maximum "sum" (column G) value across the data rows; the rows
# whose sum equals this max are flagged with Save = 1, all others get 0.
$maxVal = $null
for ($r = 2; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, 7).Value2
    if ($maxVal -eq $null -or $v -gt $maxVal) { $maxVal = $v }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -eq $maxVal) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}

Write-Host "Save column populated through row $lastRow"
